$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices, percentages) are stored as text
# so Excel does not reinterpret/round them as numbers.
$textCells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'E10', 'D11', 'E11', 'D12', 'E12', 'E13', 'D14', 'E14', 'D15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D29', 'E29', 'E30', 'D31', 'E31', 'D32', 'E32', 'E33', 'D34', 'E34', 'D35', 'E35', 'D36', 'E36', 'D37', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'E41', 'D42', 'E42', 'D43', 'E43', 'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'E48', 'E49', 'D50', 'E50', 'D51', 'E51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range('D2').Value = '96.779.51'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '3.685.16'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '243.53'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('E6').Value = '  +13.34%  '
$ws.Range('D7').Value = '670.29'
$ws.Range('E7').Value = '  +2.62%  '
$ws.Range('D8').Value = '0.428'
$ws.Range('E8').Value = '  +4.14%  '
$ws.Range('D9').Value = '1.10'
$ws.Range('E9').Value = '  +2.19%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').Value = '3.685.54'
$ws.Range('E11').Value = '  +3.36%  '
$ws.Range('D12').Value = '45.43'
$ws.Range('E12').Value = '  +4.35%  '
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('D14').Value = '6.65'
$ws.Range('E14').Value = '  +3.84%  '
$ws.Range('D15').Value = '4.373.32'
$ws.Range('D16').Value = '0.0000273'
$ws.Range('E16').Value = '  +4.93%  '
$ws.Range('D17').Value = '96.468.09'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '9.02'
$ws.Range('E18').Value = '  +16.04%  '
$ws.Range('D19').Value = '3.690.18'
$ws.Range('E19').Value = '  +3.77%  '
$ws.Range('D20').Value = '12.84'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').Value = '18.55'
$ws.Range('E21').Value = '  +3.26%  '
$ws.Range('D22').Value = '0.537'
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = '529.77'
$ws.Range('E23').Value = '  +4.24%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').Value = '3.50'
$ws.Range('E24').Value = '  +3.13%  '
$ws.Range('D25').Value = '0.0000207'
$ws.Range('E25').Value = '  +3.31%  '
$ws.Range('D26').Value = '7.07'
$ws.Range('E26').Value = '  +1.81%  '
$ws.Range('D27').Value = '103.58'
$ws.Range('E27').Value = '  +7.16%  '
$ws.Range('D28').Value = '13.07'
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('D29').Value = '0.166'
$ws.Range('E29').Value = '  +6.93%  '
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('D31').Value = '12.21'
$ws.Range('E31').Value = '  +6.29%  '
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  +1.82%  '
$ws.Range('D34').Value = '1.85'
$ws.Range('E34').Value = '  +13.05%  '
$ws.Range('D35').Value = '32.98'
$ws.Range('E35').Value = '  +5.03%  '
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = '0.590'
$ws.Range('E37').Value = '  +3.93%  '
$ws.Range('D38').Value = '624.03'
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('D39').Value = '8.87'
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('D40').Value = '42.61'
$ws.Range('E40').Value = '  +30.31%  '
$ws.Range('E41').Value = '  +6.93%  '
$ws.Range('B42').Value = 'ImmutableX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D42').Value = '1.96'
$ws.Range('E42').Value = '  +7.85%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '0.959'
$ws.Range('E43').Value = '  +5.88%  '
$ws.Range('D45').Value = '6.22'
$ws.Range('E45').Value = '  +7.40%  '
$ws.Range('D46').Value = '0.0460'
$ws.Range('E46').Value = '  +7.64%  '
$ws.Range('D47').Value = '0.424'
$ws.Range('E47').Value = '  +20.47%  '
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('D50').Value = '8.60'
$ws.Range('E50').Value = '  +3.20%  '
$ws.Range('B51').Value = 'MantraDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D51').Value = '3.58'
$ws.Range('E51').Value = '  +1.81%  '
